# Apply backfilled portfolio rows to "Daily" sheet and update "Holdings" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Daily sheet: append backfilled rows 185-188 (Feb 24-27, 2026), each
#    tagged with note "auto-append". The date column for these rows is a
#    literal text string (e.g. "2026-02-24"), NOT a real Excel date serial
#    like the rest of the column, so it is written as plain text.
# ---------------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily")

$newRows = @(
    @("2026-02-24", 571.73, 1483.58, 1696.46, 3751.77, 0.9014, "auto-append"),
    @("2026-02-25", 571.73, 1498.14, 1723.54, 3793.41, 0.9114, "auto-append"),
    @("2026-02-26", 571.73, 1500.24, 1692.01, 3763.98, 0.9043, "auto-append"),
    @("2026-02-27", 571.73, 1500.24, 1692.01, 3763.98, 0.9043, "auto-append")
)

# A scratch cell (far outside the used range) is used to stage the date as
# plain text (leading apostrophe forces text, not an auto-converted date
# serial), then copy/paste-special(values) it into column A so the target
# cell picks up the text value without inheriting a new number-format style.
# The scratch cell is fully cleared (not just its contents) afterwards so it
# leaves no residue in the sheet's used range / dimension.
$scratch = $daily.Cells.Item(500, 1)

$startRow = 185
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $scratch.Value = "'" + $row[0]
    $scratch.Copy()
    $daily.Cells.Item($r, 1).PasteSpecial(-4163)
    $scratch.Clear()

    $daily.Cells.Item($r, 2).Value = $row[1]
    $daily.Cells.Item($r, 3).Value = $row[2]
    $daily.Cells.Item($r, 4).Value = $row[3]
    $daily.Cells.Item($r, 5).Value = $row[4]
    $daily.Cells.Item($r, 6).Value = $row[5]
    $daily.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2) Holdings sheet: update Gold Futures quantity, swap NVDA->QQQ in row 4,
#    update TSLA quantity, and re-add NVDA as a new row 6.
# ---------------------------------------------------------------------------
$holdings = $wb.Worksheets.Item("Holdings")

$holdings.Cells.Item(3, 3).Value = 0.2877491817892204

$holdings.Cells.Item(4, 1).Value = "QQQ"
$holdings.Cells.Item(4, 2).Value = "Invesco QQQ Trust"
$holdings.Cells.Item(4, 3).Value = 1.7327

$holdings.Cells.Item(5, 3).Value = 1.2254

$holdings.Cells.Item(6, 1).Value = "NVDA"
$holdings.Cells.Item(6, 2).Value = "NVIDIA Corp"
$holdings.Cells.Item(6, 3).Value = 0.734
